$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.281.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.073.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.63%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.072.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.52%  "
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.486"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.38%  "
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.40%  "
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.588.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.294.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.079.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +17.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.05%  "
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.52%  "
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("E31").Value = "  +3.33%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.57%  "
$ws.Range("E34").Value = "  +5.61%  "
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  +5.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "383.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.772.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.55%  "
